# Update "想去人数" (want-to-go count) figures for the latest data refresh.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 13552
$ws.Range("F13").Value = 13564
$ws.Range("F16").Value = 8967
$ws.Range("F31").Value = 188

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 13552
$ws.Range("F13").Value = 13564
$ws.Range("F16").Value = 8967
$ws.Range("F33").Value = 188
